$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently carries a pandas-exported "index" column (col A) and an
# extra header row used for pandas-style column metadata. Remove that
# leftover numbering: drop column A (the bare row index 0,1,2) and then drop
# what is now row 1 (the original row 1, which only ever held the two
# "leftover" index/placeholder cells in B1/C1) so the data shifts up-and-left
# into a clean A1:B3 block.
$ws.Columns.Item(1).Delete()
$ws.Rows.Item(1).Delete()
